# Auto-generated edit script: updates cached market-price values
# across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match
# the latest scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 6
$ws.Range("H6").Value = 2047.8
$ws.Range("I6").Value = 2047.8
$ws.Range("K6").Value = 6143.4
$ws.Range("M6").Value = -6031.4

# row 15
$ws.Range("H15").Value = 1226.2766
$ws.Range("I15").Value = 1226.2766
$ws.Range("K15").Value = 3678.8298
$ws.Range("M15").Value = -3509.8298

# row 34
$ws.Range("H34").Value = 3797.5
$ws.Range("I34").Value = 3797.5
$ws.Range("K34").Value = 3797.5
$ws.Range("M34").Value = -3594.5

# row 36
$ws.Range("H36").Value = 3797.5
$ws.Range("I36").Value = 3797.5
$ws.Range("K36").Value = 3797.5
$ws.Range("M36").Value = -3082.5

# row 41
$ws.Range("H41").Value = 556.3333
$ws.Range("I41").Value = 556.3333
$ws.Range("K41").Value = 556.3333
$ws.Range("M41").Value = -116.3333

# row 53
$ws.Range("H53").Value = 587.875
$ws.Range("I53").Value = 657.5714
$ws.Range("J53").Value = 100
$ws.Range("K53").Value = 657.5714
$ws.Range("L53").Value = 100
$ws.Range("M53").Value = -20.57140000000004
$ws.Range("N53").Value = -1374

# row 105
$ws.Range("H105").Value = 20000
$ws.Range("I105").Value = 20000
$ws.Range("K105").Value = 20000
$ws.Range("M105").Value = -16506

# row 137
$ws.Range("H137").Value = 2484.3823
$ws.Range("I137").Value = 1512.3158
$ws.Range("K137").Value = 4536.9474
$ws.Range("M137").Value = -1986.9474

# row 138
$ws.Range("H138").Value = 4279.7896
$ws.Range("I138").Value = 3893.6843
$ws.Range("J138").Value = 4665.8945
$ws.Range("K138").Value = 11681.0529
$ws.Range("L138").Value = 13997.6835
$ws.Range("M138").Value = -6541.052899999999
$ws.Range("N138").Value = -24277.6835

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 2281.2856
$ws.Range("I2").Value = 2158.3333
$ws.Range("J2").Value = 2373.5
$ws.Range("K2").Value = 2158.3333
$ws.Range("L2").Value = 2373.5
$ws.Range("M2").Value = -2045.3333
$ws.Range("N2").Value = -2599.5

# row 9
$ws.Range("H9").Value = 30009
$ws.Range("J9").Value = 30009
$ws.Range("L9").Value = 30009
$ws.Range("N9").Value = -30349

# row 20
$ws.Range("H20").Value = 30009
$ws.Range("J20").Value = 30009
$ws.Range("L20").Value = 30009
$ws.Range("N20").Value = -30549

# row 32
$ws.Range("H32").Value = 10422.206
$ws.Range("I32").Value = 10112.102
$ws.Range("K32").Value = 10112.102
$ws.Range("M32").Value = -9825.102000000001

# row 116
$ws.Range("H116").Value = 2281.2856
$ws.Range("I116").Value = 2158.3333
$ws.Range("J116").Value = 2373.5
$ws.Range("K116").Value = 2158.3333
$ws.Range("L116").Value = 2373.5
$ws.Range("M116").Value = 135.6667000000002
$ws.Range("N116").Value = -6961.5

# row 122
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 2281.2856
$ws.Range("I3").Value = 2158.3333
$ws.Range("J3").Value = 2373.5
$ws.Range("K3").Value = 2158.3333
$ws.Range("L3").Value = 2373.5
$ws.Range("M3").Value = -2044.3333
$ws.Range("N3").Value = -2601.5

# row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("N132").ClearContents()

# row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 97.59999999999999
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# row 99
$ws.Range("H99").Value = 9198.4
$ws.Range("I99").Value = 8997.333000000001
$ws.Range("J99").Value = 9500
$ws.Range("K99").Value = 8997.333000000001
$ws.Range("L99").Value = 9500
$ws.Range("M99").Value = -7499.333000000001
$ws.Range("N99").Value = -12496

# row 126
$ws.Range("H126").Value = 9198.4
$ws.Range("I126").Value = 8997.333000000001
$ws.Range("J126").Value = 9500
$ws.Range("K126").Value = 26991.999
$ws.Range("L126").Value = 28500
$ws.Range("M126").Value = -24521.999
$ws.Range("N126").Value = -33440

# row 141
$ws.Range("H141").Value = 51339.53
$ws.Range("I141").Value = 20000
$ws.Range("J141").Value = 55518.133
$ws.Range("K141").Value = 20000
$ws.Range("L141").Value = 55518.133
$ws.Range("M141").Value = -14820
$ws.Range("N141").Value = -65878.133

$ws = $wb.Worksheets.Item("CUL")
# row 11
$ws.Range("H11").Value = 933
$ws.Range("J11").Value = 1000
$ws.Range("L11").Value = 3000
$ws.Range("N11").Value = -3280

# row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("N32").ClearContents()

# row 37
$ws.Range("H37").Value = 198999.2
$ws.Range("J37").Value = 198999.2
$ws.Range("L37").Value = 596997.6000000001
$ws.Range("N37").Value = -597221.6000000001

# row 113
$ws.Range("H113").Value = 3478.875
$ws.Range("J113").Value = 3405.1428
$ws.Range("L113").Value = 10215.4284
$ws.Range("N113").Value = -14555.4284

# row 132
$ws.Range("H132").Value = 4092.0625
$ws.Range("J132").Value = 3767.1538
$ws.Range("L132").Value = 33904.3842
$ws.Range("N132").Value = -38964.3842

$ws = $wb.Worksheets.Item("GSM")
# row 122
$ws.Range("H122").Value = 1499.5
$ws.Range("I122").Value = 1999
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5997
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -3547
$ws.Range("N122").Value = -7900

# row 126
$ws.Range("H126").Value = 13998
$ws.Range("I126").Value = 13999
$ws.Range("J126").Value = 13997
$ws.Range("K126").Value = 41997
$ws.Range("L126").Value = 41991
$ws.Range("M126").Value = -39527
$ws.Range("N126").Value = -46931

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 4999.6
$ws.Range("J22").Value = 5616.5
$ws.Range("L22").Value = 5616.5
$ws.Range("N22").Value = -6206.5

# row 27
$ws.Range("H27").Value = 4999.6
$ws.Range("J27").Value = 5616.5
$ws.Range("L27").Value = 5616.5
$ws.Range("N27").Value = -5830.5

# row 55
$ws.Range("H55").Value = 372
$ws.Range("J55").Value = 380
$ws.Range("L55").Value = 380
$ws.Range("N55").Value = -726

$ws = $wb.Worksheets.Item("WVR")
# row 132
$ws.Range("H132").Value = 2431.818
$ws.Range("I132").Value = 1680.5
$ws.Range("K132").Value = 5041.5
$ws.Range("M132").Value = -2511.5
